# JournalFY17.xlsx - "Journal" sheet corrections
#
# 1. Rows 16-30 (S.No. column A) were renumbered: the two entries that used to
#    carry S.No. 15 and 16 were removed upstream, so every subsequent S.No.
#    shifts down by 2 (17->15, 18->16, ... 31->29).
# 2. The stray, effectively-empty trailing row 31 (only a styled-but-blank
#    C31 cell) is deleted outright, shrinking the used range from A1:F31 to
#    A1:F30.
# 3. The saved cursor/selection moves to D38 (below the data, where the user
#    last clicked) and the view is scrolled down so row 16 is at the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- 1. Renumber the S.No. column for rows 16-30 (shift each value down by 2)
for ($r = 16; $r -le 30; $r++) {
    $current = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $current - 2
}

# --- 2. Remove the now-obsolete trailing row 31 entirely (shifts dimension
#        from A1:F31 down to A1:F30, matching the autoFilter range already
#        present on the sheet)
$ws.Rows(31).Delete()

# --- 3. Scroll the window so row 16 is the first visible row, then move the
#        selection to D38
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$ws.Range("D38").Select() | Out-Null
